$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-45 down to 5-46
$ws.Rows(4).Insert()

# Populate the new row's three cells
$ws.Range("A4").Value = "        新增"
$ws.Range("B4").Value = "add"
$ws.Range("C4").Value = "/addArchives"

# Style the new row with a red font color (matches the new font added to styles.xml)
$ws.Range("A4:C4").Font.Color = 255

# Update the sheet view: move the active selection to A14 (matches the
# diff's sheetView/selection change; the scrolled-down topLeftCell="A10"
# is also cleared as a side effect of reselecting)
[void]$ws.Range("A14").Select()
